$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2837.3
$ws.Range("I111").Value = 2469.9092
$ws.Range("K111").Value = 7409.7276
$ws.Range("M111").Value = -4342.7276
$ws.Range("H124").Value = 34624.445
$ws.Range("J124").Value = 34624.445
$ws.Range("L124").Value = 34624.445
$ws.Range("N124").Value = -44444.445
$ws.Range("H129").Value = 865.9524
$ws.Range("I129").Value = 532.55554
$ws.Range("J129").Value = 1116
$ws.Range("K129").Value = 1597.66662
$ws.Range("L129").Value = 3348
$ws.Range("M129").Value = 3402.33338
$ws.Range("N129").Value = -13348

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35414.914
$ws.Range("I32").Value = 13783.149
$ws.Range("J32").Value = 137084.2
$ws.Range("K32").Value = 13783.149
$ws.Range("L32").Value = 137084.2
$ws.Range("M32").Value = -13496.149
$ws.Range("N32").Value = -137658.2
$ws.Range("H61").Value = 2355.8235
$ws.Range("I61").Value = 2436.6
$ws.Range("J61").Value = 1750
$ws.Range("K61").Value = 2436.6
$ws.Range("L61").Value = 1750
$ws.Range("M61").Value = -2224.6
$ws.Range("N61").Value = -2174
$ws.Range("H74").Value = 3076.738
$ws.Range("I74").Value = 2785.2693
$ws.Range("J74").Value = 3550.375
$ws.Range("K74").Value = 2785.2693
$ws.Range("L74").Value = 3550.375
$ws.Range("M74").Value = -1911.2693
$ws.Range("N74").Value = -5298.375
$ws.Range("H77").Value = 3076.738
$ws.Range("I77").Value = 2785.2693
$ws.Range("J77").Value = 3550.375
$ws.Range("K77").Value = 13926.3465
$ws.Range("L77").Value = 17751.875
$ws.Range("M77").Value = -9558.3465
$ws.Range("N77").Value = -26487.875
$ws.Range("H136").Value = 2355.8235
$ws.Range("I136").Value = 2436.6
$ws.Range("J136").Value = 1750
$ws.Range("K136").Value = 7309.799999999999
$ws.Range("L136").Value = 5250
$ws.Range("M136").Value = -4759.799999999999
$ws.Range("N136").Value = -10350

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1344.5918
$ws.Range("I134").Value = 1332.3778
$ws.Range("J134").Value = 1482
$ws.Range("K134").Value = 3997.1334
$ws.Range("L134").Value = 4446
$ws.Range("M134").Value = -1462.1334
$ws.Range("N134").Value = -9516

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H94").Value = 7437.778
$ws.Range("I94").Value = 1950
$ws.Range("J94").Value = 9005.714
$ws.Range("K94").Value = 1950
$ws.Range("L94").Value = 9005.714
$ws.Range("M94").Value = -1499
$ws.Range("N94").Value = -9907.714
$ws.Range("H97").Value = 16615
$ws.Range("J97").Value = 16615
$ws.Range("L97").Value = 16615
$ws.Range("N97").Value = -18597
$ws.Range("H132").Value = 1705.8
$ws.Range("I132").Value = 1414.0193
$ws.Range("K132").Value = 4242.0579
$ws.Range("M132").Value = -1712.0579
$ws.Range("H133").Value = 29804.285
$ws.Range("J133").Value = 29804.285
$ws.Range("L133").Value = 29804.285
$ws.Range("N133").Value = -34864.285

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 934.913
$ws.Range("I34").Value = 276.14285
$ws.Range("J34").Value = 1223.125
$ws.Range("K34").Value = 828.4285500000001
$ws.Range("L34").Value = 3669.375
$ws.Range("M34").Value = -744.4285500000001
$ws.Range("N34").Value = -3837.375
$ws.Range("H107").Value = 838.1515000000001
$ws.Range("I107").Value = 688.5294
$ws.Range("J107").Value = 997.125
$ws.Range("K107").Value = 2065.5882
$ws.Range("L107").Value = 2991.375
$ws.Range("M107").Value = -145.5882000000001
$ws.Range("N107").Value = -6831.375
$ws.Range("H122").Value = 1275.2559
$ws.Range("I122").Value = 605.0833
$ws.Range("J122").Value = 1534.6774
$ws.Range("K122").Value = 5445.7497
$ws.Range("L122").Value = 13812.0966
$ws.Range("M122").Value = -2995.7497
$ws.Range("N122").Value = -18712.0966
$ws.Range("H131").Value = 891.1900000000001
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 891.1900000000001
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2673.57
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -12753.57

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2270.3704
$ws.Range("J43").Value = 2956.111
$ws.Range("L43").Value = 2956.111
$ws.Range("N43").Value = -3258.111
$ws.Range("H46").Value = 20097.715
$ws.Range("J46").Value = 30046
$ws.Range("L46").Value = 30046
$ws.Range("N46").Value = -30358
$ws.Range("H80").Value = 5955839
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 7146207
$ws.Range("K80").Value = 4000
$ws.Range("L80").Value = 7146207
$ws.Range("M80").Value = -3002
$ws.Range("N80").Value = -7148203
$ws.Range("H83").Value = 5955839
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 7146207
$ws.Range("K83").Value = 20000
$ws.Range("L83").Value = 35731035
$ws.Range("M83").Value = -15008
$ws.Range("N83").Value = -35741019
$ws.Range("H102").Value = 1659.2222
$ws.Range("I102").Value = 1554.4615
$ws.Range("J102").Value = 1931.6
$ws.Range("K102").Value = 1554.4615
$ws.Range("L102").Value = 1931.6
$ws.Range("M102").Value = 67.53850000000011
$ws.Range("N102").Value = -5175.6
$ws.Range("H120").Value = 36000
$ws.Range("J120").Value = 36000
$ws.Range("L120").Value = 36000
$ws.Range("N120").Value = -45676
$ws.Range("H130").Value = 49436.363
$ws.Range("J130").Value = 49436.363
$ws.Range("L130").Value = 49436.363
$ws.Range("N130").Value = -59476.363
$ws.Range("H132").Value = 2431.4285
$ws.Range("I132").Value = 2451.4866
$ws.Range("J132").Value = 2283
$ws.Range("K132").Value = 7354.459800000001
$ws.Range("L132").Value = 6849
$ws.Range("M132").Value = -4824.459800000001
$ws.Range("N132").Value = -11909

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 15772
$ws.Range("J96").Value = 15772
$ws.Range("L96").Value = 15772
$ws.Range("N96").Value = -21264
$ws.Range("H139").Value = 73788.42999999999
$ws.Range("J139").Value = 73788.42999999999
$ws.Range("L139").Value = 73788.42999999999
$ws.Range("N139").Value = -84068.42999999999

